$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 111787188
$ws.Range("B3").Value = 88915
$ws.Range("C3").Value = "Ovaliderad"
$ws.Range("D3").Value = "NT"
$ws.Range("E3").Value = 5734
$ws.Range("F3").Value = "Druvfingersvamp"
$ws.Range("G3").Value = "Ramaria botrytis"
$ws.Range("H3").Value = "(Pers.:Fr.) Bourdot"

$ws.Range("I3").Value = "'"
$ws.Range("I3").Style = "Normal"

$ws.Range("K3").Value = "'"
$ws.Range("K3").Style = "Normal"

$ws.Range("P3").Value = "Restenäs 270, Uddevalla, Boh"
$ws.Range("Q3").Value = 317495
$ws.Range("R3").Value = 6459123
$ws.Range("S3").Value = 25
$ws.Range("T3").Value = "Västra Götaland"
$ws.Range("U3").Value = "Uddevalla"
$ws.Range("V3").Value = "Bohuslän"
$ws.Range("W3").Value = "Resteröd"

$ws.Range("Y3").NumberFormat = "@"
$ws.Range("Y3").Value = "2023-08-30"
$ws.Range("Y3").Style = "Normal"

$ws.Range("Z3").Value = "14:54"

$ws.Range("AA3").NumberFormat = "@"
$ws.Range("AA3").Value = "2023-08-30"
$ws.Range("AA3").Style = "Normal"

$ws.Range("AB3").Value = "14:54"
$ws.Range("AD3").Value = $false
$ws.Range("AE3").Value = $false
$ws.Range("AG3").Value = $false
$ws.Range("AH3").Value = "Skogsmark"

$ws.Range("AT3").Value = "'"
$ws.Range("AT3").Style = "Normal"

$ws.Range("AW3").Value = "Isabell Winberg"
$ws.Range("AX3").Value = "Isabell Winberg"

$ws.Range("AY3").Value = "'"
$ws.Range("AY3").Style = "Normal"
